$wb = $excel.ActiveWorkbook

# 1. Remove the "Parameters" worksheet entirely (and its external-reference usage
#    collapses naturally once nothing references it).
$excel.DisplayAlerts = $false
[void]$wb.Worksheets("Parameters").Delete()

# 2. Work on the remaining sheet.
$ws = $wb.Worksheets("Лист1")

# D4 was stored as the text "484" - turn it into a real number, matching the
# target where every year column under this row is numeric.
$ws.Range("D4").Value = 484

# Copy the formatting of the existing "2019" column (D) across the four new
# year columns (E:H) before filling in their values, so every new cell picks
# up the same style index as its row (s=2 header row, s=12 row4, s=13 row5,
# s=9 row6) the same way Excel's own fill/copy would.
$ws.Range("D3:D6").Copy()
$ws.Range("E3:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. New year headers.
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# 4. "Number of local governments" row - constant across years.
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# 5. "Proportion ..." row.
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# 6. "Number of local governments that adopt ..." row.
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# 7. With five data columns instead of one, the data rows re-wrap shorter -
#    match the tightened row heights from the target layout.
$ws.Rows("4").RowHeight = 19.5
$ws.Rows("5").RowHeight = 44.25
$ws.Rows("6").RowHeight = 51.75

# 8. Match the saved selection state from the target file.
[void]$ws.Range("D9").Select()

$excel.DisplayAlerts = $true
